$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "explanations" tooltip row (key: explanations / Explanations / Uitleg / Erläuterungen)
$ws.Rows.Item(74).Delete()

# Remove the "clickToAddNewExplanation" tooltip row, which is now row 75
# (key: clickToAddNewExplanation / Click here to add explanation / Klik hier als u uitleg wilt toevoegen / Hier klicken, um Erläuterung hinzuzufügen)
$ws.Rows.Item(75).Delete()
